$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 74; this shifts all existing rows 74-129 down to 75-130
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with the new weekly record
$ws.Cells.Item(74, 1).Value  = 11
$ws.Cells.Item(74, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value  = "Bíobío"
$ws.Cells.Item(74, 4).Value  = 44818
$ws.Cells.Item(74, 5).Value  = 8
$ws.Cells.Item(74, 6).Value  = "Fruta"
$ws.Cells.Item(74, 7).Value  = 100108
$ws.Cells.Item(74, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(74, 9).Value  = 100108002
$ws.Cells.Item(74, 10).Value = "Mango"
$ws.Cells.Item(74, 11).Value = "Sin especificar"
$ws.Cells.Item(74, 12).Value = "Primera"
$ws.Cells.Item(74, 13).Value = 200
$ws.Cells.Item(74, 14).Value = 9000
$ws.Cells.Item(74, 15).Value = 9500
$ws.Cells.Item(74, 16).Value = 9250
$ws.Cells.Item(74, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(74, 18).Value = "Brasil"
$ws.Cells.Item(74, 19).Value = 2312
$ws.Cells.Item(74, 20).Value = 4
